$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 9.145899015680611
$ws.Range("D2").Value = 4.148400511839286
$ws.Range("E2").Value = 13.68587817337193
$ws.Range("F2").Value = 19.68404746984577
$ws.Range("G2").Value = 20.84670869758471
$ws.Range("H2").Value = 12.29785567224677
$ws.Range("I2").Value = 18.321658343957
$ws.Range("L2").Value = 10.06387030611737
$ws.Range("M2").Value = 58.85590366335928
$ws.Range("O2").Value = 17.65671855882755
$ws.Range("C3").Value = 9.246283828571988
$ws.Range("D3").Value = 4.101087588320247
$ws.Range("E3").Value = 13.45682196301882
$ws.Range("F3").Value = 19.89025507144932
$ws.Range("G3").Value = 21.18488138666099
$ws.Range("H3").Value = 12.42810655197687
$ws.Range("I3").Value = 18.64545815760614
$ws.Range("L3").Value = 9.923123591576811
$ws.Range("M3").Value = 55.35953999374686
$ws.Range("O3").Value = 17.89524458316773
$ws.Range("C4").Value = 9.313172109144361
$ws.Range("D4").Value = 4.071923752072204
$ws.Range("E4").Value = 13.32156133738739
$ws.Range("F4").Value = 20.02929859891062
$ws.Range("G4").Value = 21.41357065458133
$ws.Range("H4").Value = 12.51223198059913
$ws.Range("I4").Value = 18.85344016094728
$ws.Range("L4").Value = 9.840023427189752
$ws.Range("M4").Value = 53.08749631625405
$ws.Range("O4").Value = 18.05074953070989
$ws.Range("C5").Value = 9.34172424007698
$ws.Range("D5").Value = 4.060021484112732
$ws.Range("E5").Value = 13.26785949364632
$ws.Range("F5").Value = 20.08902125267091
$ws.Range("G5").Value = 21.51188567224809
$ws.Range("H5").Value = 12.54755272274434
$ws.Range("I5").Value = 18.94048817707499
$ws.Range("L5").Value = 9.807029146828707
$ws.Range("M5").Value = 52.13026893778299
$ws.Range("O5").Value = 18.11636467665917
$ws.Range("C6").Value = 9.346542759417176
$ws.Range("D6").Value = 4.058044366339313
$ws.Range("E6").Value = 13.25902982790774
$ws.Range("F6").Value = 20.09912118742555
$ws.Range("G6").Value = 21.52851532004852
$ws.Range("H6").Value = 12.55348031887444
$ws.Range("I6").Value = 18.95508064800817
$ws.Range("L6").Value = 9.801603998730732
$ws.Range("M6").Value = 51.96943409442011
$ws.Range("O6").Value = 18.12739480894212
$ws.Range("C7").Value = 9.313551964551179
$ws.Range("D7").Value = 4.071763291517155
$ws.Range("E7").Value = 13.32083127046436
$ws.Range("F7").Value = 20.0300917286892
$ws.Range("G7").Value = 21.41487604584611
$ws.Range("H7").Value = 12.5127041279671
$ws.Range("I7").Value = 18.85460484839371
$ws.Range("L7").Value = 9.839574887957761
$ws.Range("M7").Value = 53.074713512437
$ws.Range("O7").Value = 18.05162538530502
$ws.Range("C8").Value = 9.179405691058532
$ws.Range("D8").Value = 4.132116986411867
$ws.Range("E8").Value = 13.60582248875627
$ws.Range("F8").Value = 19.75252841025154
$ws.Range("G8").Value = 20.95883678250124
$ws.Range("H8").Value = 12.34190169859194
$ws.Range("I8").Value = 18.43139446752656
$ws.Range("L8").Value = 10.01467349232377
$ws.Range("M8").Value = 57.67642492181633
$ws.Range("O8").Value = 17.73706580127178
$ws.Range("C9").Value = 8.959249235043549
$ws.Range("D9").Value = 4.24910639290106
$ws.Range("E9").Value = 14.2044363412187
$ws.Range("F9").Value = 19.30982799793007
$ws.Range("G9").Value = 20.23956712528417
$ws.Range("H9").Value = 12.04011839667758
$ws.Range("I9").Value = 17.67479362994258
$ws.Range("L9").Value = 10.38277865899514
$ws.Range("M9").Value = 65.702954612356
$ws.Range("O9").Value = 17.19342815316688
$ws.Range("C10").Value = 8.825464005300313
$ws.Range("D10").Value = 4.333623478400617
$ws.Range("E10").Value = 14.66447121800548
$ws.Range("F10").Value = 19.05067911349196
$ws.Range("G10").Value = 19.82928132666515
$ws.Range("H10").Value = 11.83893819369574
$ws.Range("I10").Value = 17.16449974536286
$ws.Range("L10").Value = 10.66613220193693
$ws.Range("M10").Value = 70.99059548146445
$ws.Range("O10").Value = 16.84072006965356
$ws.Range("C11").Value = 8.771098756198125
$ws.Range("D11").Value = 4.371638552450988
$ws.Range("E11").Value = 14.87730655537552
$ws.Range("F11").Value = 18.94806454407263
$ws.Range("G11").Value = 19.67089011756925
$ws.Range("H11").Value = 11.75195456684492
$ws.Range("I11").Value = 16.9424853495399
$ws.Range("L11").Value = 10.79737791041067
$ws.Range("M11").Value = 73.26415495466185
$ws.Range("O11").Value = 16.6908978688387
$ws.Range("C12").Value = 8.751483668288399
$ws.Range("D12").Value = 4.385961919048323
$ws.Range("E12").Value = 14.95834609169045
$ws.Range("F12").Value = 18.91148171829154
$ws.Range("G12").Value = 19.61519621317324
$ws.Range("H12").Value = 11.71967560949639
$ws.Range("I12").Value = 16.85989339022878
$ws.Range("L12").Value = 10.84737676562338
$ws.Range("M12").Value = 74.10619888971421
$ws.Range("O12").Value = 16.63573697263057
$ws.Range("C13").Value = 8.755664279578779
$ws.Range("D13").Value = 4.38288048786799
$ws.Range("E13").Value = 14.94087412506914
$ws.Range("F13").Value = 18.91925803219631
$ws.Range("G13").Value = 19.62699669899287
$ws.Range("H13").Value = 11.72659798233197
$ws.Range("I13").Value = 16.87761479290711
$ws.Range("L13").Value = 10.83659592211951
$ws.Range("M13").Value = 73.92568938670122
$ws.Range("O13").Value = 16.64754612043265
$ws.Range("C14").Value = 8.769465309381578
$ws.Range("D14").Value = 4.372818445541536
$ws.Range("E14").Value = 14.88396522543873
$ws.Range("F14").Value = 18.94500880143811
$ws.Range("G14").Value = 19.66622096722767
$ws.Range("H14").Value = 11.74928568622918
$ws.Range("I14").Value = 16.93566064393773
$ws.Range("L14").Value = 10.80148557133732
$ws.Range("M14").Value = 73.33380891113067
$ws.Range("O14").Value = 16.68632794221007
$ws.Range("C15").Value = 8.77804656268926
$ws.Range("D15").Value = 4.366645472902155
$ws.Range("E15").Value = 14.84916262643076
$ws.Range("F15").Value = 18.96108053302786
$ws.Range("G15").Value = 19.69081174786861
$ws.Range("H15").Value = 11.76326873483887
$ws.Range("I15").Value = 16.97140893917566
$ws.Range("L15").Value = 10.78001727695824
$ws.Range("M15").Value = 72.96880506659008
$ws.Range("O15").Value = 16.71028924591452
$ws.Range("C16").Value = 8.829150786474882
$ws.Range("D16").Value = 4.331129546497907
$ws.Range("E16").Value = 14.65062767254034
$ws.Range("F16").Value = 19.05769926411377
$ws.Range("G16").Value = 19.84021999085018
$ws.Range("H16").Value = 11.844714675669
$ws.Range("I16").Value = 17.17921460825352
$ws.Range("L16").Value = 10.65759889227263
$ws.Range("M16").Value = 70.8393644773236
$ws.Range("O16").Value = 16.8507289523717
$ws.Range("C17").Value = 8.862191328874101
$ws.Range("D17").Value = 4.30922377273783
$ws.Range("E17").Value = 14.52969564068276
$ws.Range("F17").Value = 19.12093928326897
$ws.Range("G17").Value = 19.9392698017773
$ws.Range("H17").Value = 11.89584584812341
$ws.Range("I17").Value = 17.30930547731734
$ws.Range("L17").Value = 10.58307250759352
$ws.Range("M17").Value = 69.49929049673592
$ws.Range("O17").Value = 16.93963726909888
$ws.Range("C18").Value = 8.881804258919246
$ws.Range("D18").Value = 4.296584046275397
$ws.Range("E18").Value = 14.46047810220089
$ws.Range("F18").Value = 19.15874609185537
$ws.Range("G18").Value = 19.99888595468576
$ws.Range("H18").Value = 11.92568142952059
$ws.Range("I18").Value = 17.38508127314905
$ws.Range("L18").Value = 10.54043003887487
$ws.Range("M18").Value = 68.71610384660585
$ws.Range("O18").Value = 16.99177230109662
$ws.Range("C19").Value = 8.88854831110303
$ws.Range("D19").Value = 4.292297874035548
$ws.Range("E19").Value = 14.43710266359575
$ws.Range("F19").Value = 19.17179053719121
$ws.Range("G19").Value = 20.01951857540693
$ws.Range("H19").Value = 11.93585622992027
$ws.Range("I19").Value = 17.41090028186468
$ws.Range("L19").Value = 10.52603154649117
$ws.Range("M19").Value = 68.44879738832633
$ws.Range("O19").Value = 17.00959428086438
$ws.Range("C20").Value = 8.858610820671792
$ws.Range("D20").Value = 4.311559903876983
$ws.Range("E20").Value = 14.54253444673071
$ws.Range("F20").Value = 19.11405841710943
$ws.Range("G20").Value = 19.92845039806342
$ws.Range("H20").Value = 11.89035867391439
$ws.Range("I20").Value = 17.2953584913451
$ws.Range("L20").Value = 10.590983157167
$ws.Range("M20").Value = 69.64322790000821
$ws.Range("O20").Value = 16.93006928840307
$ws.Range("C21").Value = 8.765384924206829
$ws.Range("D21").Value = 4.375775946060364
$ws.Range("E21").Value = 14.90066925048689
$ws.Range("F21").Value = 18.93738279727405
$ws.Range("G21").Value = 19.65458174428172
$ws.Range("H21").Value = 11.74260378401785
$ws.Range("I21").Value = 16.91857078918358
$ws.Range("L21").Value = 10.81179051999284
$ws.Range("M21").Value = 73.50817109016407
$ws.Range("O21").Value = 16.67489369917914
$ws.Range("C22").Value = 8.710140114885839
$ws.Range("D22").Value = 4.417320139892223
$ws.Range("E22").Value = 15.13728570531422
$ws.Range("F22").Value = 18.83521336918084
$ws.Range("G22").Value = 19.5006776370648
$ws.Range("H22").Value = 11.64988761970785
$ws.Range("I22").Value = 16.68095901764252
$ws.Range("L22").Value = 10.95782646458408
$ws.Range("M22").Value = 75.92398658534258
$ws.Range("O22").Value = 16.51731662133108
$ws.Range("C23").Value = 8.739092197185098
$ws.Range("D23").Value = 4.395189275883757
$ws.Range("E23").Value = 15.01078733674188
$ws.Range("F23").Value = 18.88849998915659
$ws.Range("G23").Value = 19.58044912576278
$ws.Range("H23").Value = 11.69901694570854
$ws.Range("I23").Value = 16.80697726778801
$ws.Range("L23").Value = 10.87973884669911
$ws.Range("M23").Value = 74.64467714396977
$ws.Range("O23").Value = 16.60056134043692
$ws.Range("C24").Value = 8.860227647962441
$ws.Range("D24").Value = 4.310503880901971
$ws.Range("E24").Value = 14.53672905821808
$ws.Range("F24").Value = 19.11716474847106
$ws.Range("G24").Value = 19.93333354746756
$ws.Range("H24").Value = 11.89283805673774
$ws.Range("I24").Value = 17.30166085803143
$ws.Range("L24").Value = 10.58740611494657
$ws.Range("M24").Value = 69.57819354903778
$ws.Range("O24").Value = 16.93439179738543
$ws.Range("C25").Value = 9.014034145972282
$ws.Range("D25").Value = 4.217667030888695
$ws.Range("E25").Value = 14.0386526210789
$ws.Range("F25").Value = 19.41828935197196
$ws.Range("G25").Value = 20.41425558608268
$ws.Range("H25").Value = 12.11817678193249
$ws.Range("I25").Value = 17.87154129904593
$ws.Range("L25").Value = 10.28077159938315
$ws.Range("M25").Value = 63.63863421308407
$ws.Range("O25").Value = 17.33245335165062
